# Insert a new data row at row 351 (pushes existing rows 351-426 down to 352-427)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A351").EntireRow.Insert()

# Populate the newly-inserted row 351 with the new record's data.
$ws.Range("A351").Value = 3
$ws.Range("B351").Value = "Femacal de La Calera"
$ws.Range("C351").Value = "Coquimbo"
$ws.Range("D351").Value = 44798
$ws.Range("E351").Value = 5
$ws.Range("F351").Value = 100112031
$ws.Range("G351").Value = "Poroto verde"
$ws.Range("H351").Value = "Magnum"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 45
$ws.Range("K351").Value = 33000
$ws.Range("L351").Value = 33000
$ws.Range("M351").Value = 33000
$ws.Range("N351").Value = "$/malla 25 kilos"
$ws.Range("O351").Value = "Región de Arica y Parinacota"
$ws.Range("P351").Value = 1320
$ws.Range("Q351").Value = 25
$ws.Range("R351").Value = "Hortaliza"
